$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.703.23'
$ws.Range("E2").Value = '  +0.68%  '
$ws.Range("D3").Value = '3.691.40'
$ws.Range("E3").Value = '  +0.51%  '
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '671.22'
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.52'
$ws.Range("E6").Value = '  +1.76%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +1.53%  '
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.11'
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.442'
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '33.13'
$ws.Range("E13").Value = '  +2.55%  '
$ws.Range("D14").Value = '3.626.60'
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").Value = '69.699.77'
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("E16").Value = '  +2.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '16.17'
$ws.Range("E17").Value = '  +0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.48'
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '471.27'
$ws.Range("E19").Value = '  +0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.76'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("E21").Value = '  -0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '79.88'
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '3.838.53'
$ws.Range("E24").Value = '  +5.00%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.98'
$ws.Range("E26").Value = '  +0.50%  '
$ws.Range("E27").Value = '  +0.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.69'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.71'
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.02'
$ws.Range("E30").Value = '  +1.71%  '
$ws.Range("E31").Value = '  +4.30%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.85'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("E34").Value = '  -1.81%  '
$ws.Range("D35").Value = '3.689.86'
$ws.Range("E35").Value = '  +0.65%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.51'
$ws.Range("E36").Value = '  +4.61%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.12'
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("E38").Value = '  -0.07%  '
$ws.Range("E39").Value = '  +2.01%  '
$ws.Range("E40").Value = '  +0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '177.09'
$ws.Range("E41").Value = '  +1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0910'
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("E43").Value = '  -0.35%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '47.03'
$ws.Range("E45").Value = '  +2.11%  '
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '27.55'
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("E48").Value = '  -1.21%  '
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '366.26'
$ws.Range("E51").Value = '  +1.79%  '
